# Update the cryptos price/volume snapshot (GitHub Actions refresh).
# Column D (Price) and E (Volume(1h)) get refreshed numbers for most rows;
# rows 13/14 and 47/48 also swap rank order (names, links, price, volume).
# For D-column values that look like plain numbers, force the cell to
# Text format first so Excel stores the exact digits/trailing zeros
# instead of silently re-parsing them as a floating point number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.295.06'
$ws.Range('E2').Value = '  -4.94%  '
$ws.Range('D3').Value = '1.840.45'
$ws.Range('E3').Value = '  -4.67%  '
$ws.Range('E4').Value = '  -0.75%  '
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('E6').Value = '  -0.53%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4595'
$ws.Range('E7').Value = '  -4.56%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3863'
$ws.Range('E8').Value = '  -5.87%  '
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07865'
$ws.Range('E10').Value = '  -3.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.9636'
$ws.Range('E11').Value = '  -4.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.95'
$ws.Range('E12').Value = '  -7.03%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.723'
$ws.Range('E13').Value = '  -5.43%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.793.58'
$ws.Range('E14').Value = '  -7.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.924'
$ws.Range('E15').Value = '  -4.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06834'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '87.00'
$ws.Range('E18').Value = '  -4.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000009953'
$ws.Range('E19').Value = '  -3.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.92'
$ws.Range('E20').Value = '  -4.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.49%  '
$ws.Range('D22').Value = '28.329.45'
$ws.Range('E22').Value = '  -4.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.342'
$ws.Range('E23').Value = '  -4.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.96'
$ws.Range('E24').Value = '  -7.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.141'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').Value = '2.039.49'
$ws.Range('E26').Value = '  -6.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.61'
$ws.Range('E27').Value = '  -1.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.19'
$ws.Range('E28').Value = '  -3.91%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.755'
$ws.Range('E29').Value = '  -13.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.989'
$ws.Range('E30').Value = '  -4.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '117.17'
$ws.Range('E31').Value = '  -3.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9448'
$ws.Range('E32').Value = '  -5.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09294'
$ws.Range('E33').Value = '  -3.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.286'
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.442'
$ws.Range('E35').Value = '  -2.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.327'
$ws.Range('E36').Value = '  -5.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06022'
$ws.Range('E37').Value = '  -8.07%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02153'
$ws.Range('E38').Value = '  -5.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.146'
$ws.Range('E39').Value = '  -4.67%  '
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.636'
$ws.Range('E41').Value = '  -3.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5633'
$ws.Range('E42').Value = '  -5.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.01'
$ws.Range('E43').Value = '  -6.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1782'
$ws.Range('E44').Value = '  -3.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.236'
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.279'
$ws.Range('E46').Value = '  -7.91%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '11.58'
$ws.Range('E47').Value = '  -5.54%  '
$ws.Range('B48').Value = 'Decentraland'
$ws.Range('C48').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5288'
$ws.Range('E48').Value = '  -4.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07043'
$ws.Range('E49').Value = '  -5.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.837'
$ws.Range('E50').Value = '  -7.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '112.41'
$ws.Range('E51').Value = '  -3.65%  '
